$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (ALC)
$ws.Cells.Item(18, 8).Value = 378
$ws.Cells.Item(18, 9).Value = 364.44446
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 364.44446
$ws.Cells.Item(18, 12).Value = 500
$ws.Cells.Item(18, 13).Value = -80.44445999999999
$ws.Cells.Item(18, 14).Value = -1068

# Row 40 (ALC)
$ws.Cells.Item(40, 8).Value = 2101.2727
$ws.Cells.Item(40, 10).Value = 2472.3333
$ws.Cells.Item(40, 12).Value = 2472.3333
$ws.Cells.Item(40, 14).Value = -2822.3333

# Row 49 (ALC)
$ws.Cells.Item(49, 8).Value = 5349.75
$ws.Cells.Item(49, 10).Value = 6633
$ws.Cells.Item(49, 12).Value = 19899
$ws.Cells.Item(49, 14).Value = -20171

# Row 106 (ALC)
$ws.Cells.Item(106, 8).Value = 57358.332
$ws.Cells.Item(106, 9).Value = 66130.10000000001
$ws.Cells.Item(106, 11).Value = 66130.10000000001
$ws.Cells.Item(106, 13).Value = -65499.10000000001

# Row 111 (ALC)
$ws.Cells.Item(111, 8).Value = 2150.625
$ws.Cells.Item(111, 9).Value = 1695.6
$ws.Cells.Item(111, 10).Value = 2909
$ws.Cells.Item(111, 11).Value = 5086.799999999999
$ws.Cells.Item(111, 12).Value = 8727
$ws.Cells.Item(111, 13).Value = -2019.799999999999
$ws.Cells.Item(111, 14).Value = -14861

# Row 116 (ALC)
$ws.Cells.Item(116, 8).Value = 5632.3213
$ws.Cells.Item(116, 9).Value = 4755.1113
$ws.Cells.Item(116, 10).Value = 6047.8423
$ws.Cells.Item(116, 11).Value = 4755.1113
$ws.Cells.Item(116, 12).Value = 6047.8423
$ws.Cells.Item(116, 13).Value = -1313.1113
$ws.Cells.Item(116, 14).Value = -12931.8423

# Row 118 (ALC)
$ws.Cells.Item(118, 8).Value = 479.18182
$ws.Cells.Item(118, 9).Value = 497.6
$ws.Cells.Item(118, 11).Value = 1492.8
$ws.Cells.Item(118, 13).Value = 164.1999999999998

# Row 123 (ALC)
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).Value = ""

# Row 136 (ALC)
$ws.Cells.Item(136, 8).Value = 79914.664
$ws.Cells.Item(136, 10).Value = 79914.664
$ws.Cells.Item(136, 12).Value = 79914.664
$ws.Cells.Item(136, 14).Value = -90114.664

# Row 137 (ALC)
$ws.Cells.Item(137, 8).Value = 2092.7144
$ws.Cells.Item(137, 9).Value = 2251.8667
$ws.Cells.Item(137, 11).Value = 6755.6001
$ws.Cells.Item(137, 13).Value = -4205.6001

# Row 138 (ALC)
$ws.Cells.Item(138, 8).Value = 2213.7334
$ws.Cells.Item(138, 10).Value = 2307.5715
$ws.Cells.Item(138, 12).Value = 6922.7145
$ws.Cells.Item(138, 14).Value = -17202.7145

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Cells.Item(61, 8).Value = 4566.0713
$ws.Cells.Item(61, 9).Value = 3725.25
$ws.Cells.Item(61, 11).Value = 3725.25
$ws.Cells.Item(61, 13).Value = -3513.25

# Row 74 (ARM)
$ws.Cells.Item(74, 8).Value = 3785.0667
$ws.Cells.Item(74, 9).Value = 3660.75
$ws.Cells.Item(74, 11).Value = 3660.75
$ws.Cells.Item(74, 13).Value = -2786.75

# Row 77 (ARM)
$ws.Cells.Item(77, 8).Value = 3785.0667
$ws.Cells.Item(77, 9).Value = 3660.75
$ws.Cells.Item(77, 11).Value = 18303.75
$ws.Cells.Item(77, 13).Value = -13935.75

# Row 119 (ARM)
$ws.Cells.Item(119, 8).Value = 25000
$ws.Cells.Item(119, 10).Value = 25000
$ws.Cells.Item(119, 12).Value = 25000
$ws.Cells.Item(119, 14).Value = -34676

# Row 132 (ARM)
$ws.Cells.Item(132, 8).Value = 6034
$ws.Cells.Item(132, 9).Value = 3812.805
$ws.Cells.Item(132, 10).Value = 15140.9
$ws.Cells.Item(132, 11).Value = 11438.415
$ws.Cells.Item(132, 12).Value = 45422.7
$ws.Cells.Item(132, 13).Value = -8908.414999999999
$ws.Cells.Item(132, 14).Value = -50482.7

# Row 136 (ARM)
$ws.Cells.Item(136, 8).Value = 4566.0713
$ws.Cells.Item(136, 9).Value = 3725.25
$ws.Cells.Item(136, 11).Value = 11175.75
$ws.Cells.Item(136, 13).Value = -8625.75

$ws = $wb.Worksheets.Item("BSM")
# Row 59 (BSM)
$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).Value = ""

# Row 86 (BSM)
$ws.Cells.Item(86, 8).Value = 3060.5
$ws.Cells.Item(86, 9).Value = 3729.6
$ws.Cells.Item(86, 10).Value = 2391.4
$ws.Cells.Item(86, 11).Value = 3729.6
$ws.Cells.Item(86, 12).Value = 2391.4
$ws.Cells.Item(86, 13).Value = -2606.6
$ws.Cells.Item(86, 14).Value = -4637.4

# Row 89 (BSM)
$ws.Cells.Item(89, 8).Value = 3060.5
$ws.Cells.Item(89, 9).Value = 3729.6
$ws.Cells.Item(89, 10).Value = 2391.4
$ws.Cells.Item(89, 11).Value = 18648
$ws.Cells.Item(89, 12).Value = 11957
$ws.Cells.Item(89, 13).Value = -13032
$ws.Cells.Item(89, 14).Value = -23189

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Cells.Item(16, 8).Value = 2478.6
$ws.Cells.Item(16, 9).Value = 2999
$ws.Cells.Item(16, 10).Value = 2348.5
$ws.Cells.Item(16, 11).Value = 2999
$ws.Cells.Item(16, 12).Value = 2348.5
$ws.Cells.Item(16, 13).Value = -2712
$ws.Cells.Item(16, 14).Value = -2922.5

# Row 35 (CRP)
$ws.Cells.Item(35, 8).Value = 527.6
$ws.Cells.Item(35, 9).Value = 434.5
$ws.Cells.Item(35, 10).Value = 900
$ws.Cells.Item(35, 11).Value = 434.5
$ws.Cells.Item(35, 12).Value = 900
$ws.Cells.Item(35, 13).Value = -140.5
$ws.Cells.Item(35, 14).Value = -1488

# Row 41 (CRP)
$ws.Cells.Item(41, 8).Value = 4183.1665
$ws.Cells.Item(41, 9).Value = 4183.1665
$ws.Cells.Item(41, 11).Value = 4183.1665
$ws.Cells.Item(41, 13).Value = -3755.1665

# Row 58 (CRP)
$ws.Cells.Item(58, 8).Value = 2161
$ws.Cells.Item(58, 9).Value = 2134.4546
$ws.Cells.Item(58, 11).Value = 2134.4546
$ws.Cells.Item(58, 13).Value = -1931.4546

# Row 113 (CRP)
$ws.Cells.Item(113, 8).Value = 2478.6
$ws.Cells.Item(113, 9).Value = 2999
$ws.Cells.Item(113, 10).Value = 2348.5
$ws.Cells.Item(113, 11).Value = 2999
$ws.Cells.Item(113, 12).Value = 2348.5
$ws.Cells.Item(113, 13).Value = -829
$ws.Cells.Item(113, 14).Value = -6688.5

# Row 122 (CRP)
$ws.Cells.Item(122, 8).Value = 2002.0416
$ws.Cells.Item(122, 9).Value = 1938.5834
$ws.Cells.Item(122, 11).Value = 5815.7502
$ws.Cells.Item(122, 13).Value = -3365.7502

# Row 136 (CRP)
$ws.Cells.Item(136, 8).Value = 2161
$ws.Cells.Item(136, 9).Value = 2134.4546
$ws.Cells.Item(136, 11).Value = 6403.3638
$ws.Cells.Item(136, 13).Value = -3853.3638

$ws = $wb.Worksheets.Item("CUL")
# Row 59 (CUL)
$ws.Cells.Item(59, 8).Value = 1806.125
$ws.Cells.Item(59, 9).Value = 1749
$ws.Cells.Item(59, 10).Value = 1809.9333
$ws.Cells.Item(59, 11).Value = 5247
$ws.Cells.Item(59, 12).Value = 5429.7999
$ws.Cells.Item(59, 13).Value = -4707
$ws.Cells.Item(59, 14).Value = -6509.7999

$ws = $wb.Worksheets.Item("GSM")
# Row 126 (GSM)
$ws.Cells.Item(126, 8).Value = 9370.125
$ws.Cells.Item(126, 9).Value = 9192.4
$ws.Cells.Item(126, 11).Value = 27577.2
$ws.Cells.Item(126, 13).Value = -25107.2

# Row 132 (GSM)
$ws.Cells.Item(132, 8).Value = 1160.4286
$ws.Cells.Item(132, 9).Value = 1103.8334
$ws.Cells.Item(132, 11).Value = 3311.5002
$ws.Cells.Item(132, 13).Value = -781.5001999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Cells.Item(61, 8).Value = 6166.25
$ws.Cells.Item(61, 9).Value = 5444
$ws.Cells.Item(61, 10).Value = 6599.6
$ws.Cells.Item(61, 11).Value = 5444
$ws.Cells.Item(61, 12).Value = 6599.6
$ws.Cells.Item(61, 13).Value = -5242
$ws.Cells.Item(61, 14).Value = -7003.6

# Row 82 (LTW)
$ws.Cells.Item(82, 8).Value = 1961.8125
$ws.Cells.Item(82, 10).Value = 1879.4
$ws.Cells.Item(82, 12).Value = 1879.4
$ws.Cells.Item(82, 14).Value = -2601.4

# Row 85 (LTW)
$ws.Cells.Item(85, 8).Value = 1961.8125
$ws.Cells.Item(85, 10).Value = 1879.4
$ws.Cells.Item(85, 12).Value = 1879.4
$ws.Cells.Item(85, 14).Value = -4375.4

# Row 86 (LTW)
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).Value = ""

# Row 89 (LTW)
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).Value = ""

# Row 113 (LTW)
$ws.Cells.Item(113, 8).Value = 6166.25
$ws.Cells.Item(113, 9).Value = 5444
$ws.Cells.Item(113, 10).Value = 6599.6
$ws.Cells.Item(113, 11).Value = 5444
$ws.Cells.Item(113, 12).Value = 6599.6
$ws.Cells.Item(113, 13).Value = -3274
$ws.Cells.Item(113, 14).Value = -10939.6

# Row 122 (LTW)
$ws.Cells.Item(122, 8).Value = 6755.1377
$ws.Cells.Item(122, 9).Value = 5497.1
$ws.Cells.Item(122, 11).Value = 16491.3
$ws.Cells.Item(122, 13).Value = -14041.3

$ws = $wb.Worksheets.Item("WVR")
# Row 33 (WVR)
$ws.Cells.Item(33, 8).Value = 30019
$ws.Cells.Item(33, 9).Value = 30019
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 30019
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = -29769
$ws.Cells.Item(33, 14).Value = ""

# Row 36 (WVR)
$ws.Cells.Item(36, 8).Value = 30019
$ws.Cells.Item(36, 9).Value = 30019
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 30019
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -29769
$ws.Cells.Item(36, 14).Value = ""

# Row 41 (WVR)
$ws.Cells.Item(41, 8).Value = 9601.308000000001
$ws.Cells.Item(41, 10).Value = 8552
$ws.Cells.Item(41, 12).Value = 8552
$ws.Cells.Item(41, 14).Value = -9332

# Row 44 (WVR)
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).Value = ""

# Row 62 (WVR)
$ws.Cells.Item(62, 8).Value = 68936.12
$ws.Cells.Item(62, 9).Value = 103448.27
$ws.Cells.Item(62, 10).Value = 5663.8335
$ws.Cells.Item(62, 11).Value = 103448.27
$ws.Cells.Item(62, 12).Value = 5663.8335
$ws.Cells.Item(62, 13).Value = -102824.27
$ws.Cells.Item(62, 14).Value = -6911.8335

# Row 65 (WVR)
$ws.Cells.Item(65, 8).Value = 68936.12
$ws.Cells.Item(65, 9).Value = 103448.27
$ws.Cells.Item(65, 10).Value = 5663.8335
$ws.Cells.Item(65, 11).Value = 517241.35
$ws.Cells.Item(65, 12).Value = 28319.1675
$ws.Cells.Item(65, 13).Value = -514121.35
$ws.Cells.Item(65, 14).Value = -34559.1675

# Row 113 (WVR)
$ws.Cells.Item(113, 8).Value = 1814.6666
$ws.Cells.Item(113, 9).Value = 1466.6923
$ws.Cells.Item(113, 10).Value = 2719.4
$ws.Cells.Item(113, 11).Value = 4400.0769
$ws.Cells.Item(113, 12).Value = 8158.200000000001
$ws.Cells.Item(113, 13).Value = -2230.0769
$ws.Cells.Item(113, 14).Value = -12498.2

# Row 132 (WVR)
$ws.Cells.Item(132, 8).Value = 1461.75
$ws.Cells.Item(132, 9).Value = 1565.0667
$ws.Cells.Item(132, 10).Value = 1151.8
$ws.Cells.Item(132, 11).Value = 4695.2001
$ws.Cells.Item(132, 12).Value = 3455.4
$ws.Cells.Item(132, 13).Value = -2165.2001
$ws.Cells.Item(132, 14).Value = -8515.4

# Row 136 (WVR)
$ws.Cells.Item(136, 8).Value = 3030.9546
$ws.Cells.Item(136, 10).Value = 4125
$ws.Cells.Item(136, 12).Value = 12375
$ws.Cells.Item(136, 14).Value = -17475
